$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

$ws.Range("E2").Value = "RMA-REX3-001"
$ws.Range("F2").Value = "RMA-REX3-1-1"
$ws.Range("J2").Value = "a7s5f000000xK6fAAE"

$ws.Range("E3").Value = "RMA-REX3-002"
$ws.Range("F3").Value = "RMA-REX3-1-2"
$ws.Range("J3").Value = "a7s5f000000xK6gAAE"

$ws.Range("E4").Value = "RMA-REX3-003"
$ws.Range("F4").Value = "RMA-REX3-1-3"
$ws.Range("J4").Value = "a7s5f000000xK6hAAE"

$ws.Activate()
